$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force the price/volume columns to remain text so numeric-looking strings
# (e.g. '245.68', '1.0000') are not auto-converted to numbers by Excel.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '29.198.60'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '1.848.42'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '245.68'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').Value = '0.7034'
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('D7').Value = '0.9998'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '0.07782'
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '0.3071'
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('D10').Value = '23.63'
$ws.Range('E10').Value = '  -0.73%  '
$ws.Range('D11').Value = '0.07813'
$ws.Range('D12').Value = '93.26'
$ws.Range('E12').Value = '  +1.17%  '
$ws.Range('D13').Value = '5.146'
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('D14').Value = '1.848.57'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').Value = '0.6873'
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').Value = '6.601'
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('D17').Value = '0.000008346'
$ws.Range('E17').Value = '  -1.06%  '
$ws.Range('D18').Value = '29.195.50'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').Value = '242.28'
$ws.Range('E19').Value = '  -2.95%  '
$ws.Range('D20').Value = '2.092.54'
$ws.Range('E20').Value = '  -0.68%  '
$ws.Range('D21').Value = '12.74'
$ws.Range('E21').Value = '  -0.62%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').Value = '7.520'
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').Value = '1.0000'
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('E25').Value = '  +1.22%  '
$ws.Range('D26').Value = '159.34'
$ws.Range('E26').Value = '  -1.15%  '
$ws.Range('D27').Value = '8.846'
$ws.Range('E27').Value = '  -0.25%  '
$ws.Range('D29').Value = '1.537'
$ws.Range('E29').Value = '  -1.48%  '
$ws.Range('D30').Value = '4.226'
$ws.Range('E30').Value = '  -0.53%  '
$ws.Range('E31').Value = '  -0.62%  '
$ws.Range('D32').Value = '1.198'
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('D33').Value = '0.05123'
$ws.Range('E33').Value = '  -1.54%  '
$ws.Range('D34').Value = '0.7903'
$ws.Range('E34').Value = '  +3.94%  '
$ws.Range('D35').Value = '1.894'
$ws.Range('E35').Value = '  +2.72%  '
$ws.Range('E36').Value = '  -1.90%  '
$ws.Range('D37').Value = '2.694'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('D38').Value = '1.317.36'
$ws.Range('E38').Value = '  +7.94%  '
$ws.Range('D39').Value = '0.01870'
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').Value = '2.713'
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('D41').Value = '0.9545'
$ws.Range('E41').Value = '  +6.14%  '
$ws.Range('D42').Value = '6.082'
$ws.Range('E42').Value = '  +10.43%  '
$ws.Range('D43').Value = '107.02'
$ws.Range('E43').Value = '  -2.61%  '
$ws.Range('D44').Value = '0.9999'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').Value = '9.693'
$ws.Range('E45').Value = '  +1.25%  '
$ws.Range('D46').Value = '1.992.77'
$ws.Range('E46').Value = '  -0.69%  '
$ws.Range('E47').Value = '  +0.09%  '
$ws.Range('D48').Value = '64.41'
$ws.Range('E48').Value = '  -1.76%  '
$ws.Range('D49').Value = '1.765'
$ws.Range('E49').Value = '  +0.61%  '
$ws.Range('D50').Value = '6.989'
$ws.Range('E50').Value = '  -0.55%  '
$ws.Range('E51').Value = '  -0.73%  '

# Restore the default (unstyled) look so no stray style index is left behind.
$ws.Range('D2:E51').Style = 'Normal'
